# Horarios Linea 141 - actualizacion de datos (scrape 04:43:39)
$wb = $excel.ActiveWorkbook

$newTime = "04:43:39"

# ---------------------------------------------------------------------------
# Hoja 1: LP1912  (10 filas -> 14 filas)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 14"

$rows1 = @(
    @($newTime, "04:45", "215A_EL PATO", 2, "LP1912"),
    @($newTime, "04:53", "11_ETCHEVERRY", 10, "LP1912"),
    @($newTime, "05:16", "17_ROMERO", 33, "LP1912"),
    @($newTime, "05:22", "23_HERNANDEZ", 39, "LP1912"),
    @($newTime, "05:34", "215B_EL PATO", 51, "LP1912"),
    @($newTime, "05:46", "15_ABASTO", 63, "LP1912"),
    @($newTime, "05:54", "10_OLMOS", 71, "LP1912"),
    @($newTime, "06:04", "16_SANTA ANA", 81, "LP1912"),
    @($newTime, "06:11", "215A_EL PATO", 88, "LP1912"),
    @($newTime, "06:14", "225_HARAS DEL SUR", 91, "LP1912"),
    @($newTime, "06:21", "26_HERNANDEZ", 98, "LP1912"),
    @($newTime, "06:27", "23_HERNANDEZ", 104, "LP1912"),
    @($newTime, "06:29", "86_EST CHICA-ESC AGRARIA", 106, "LP1912"),
    @($newTime, "06:31", "16_SANTA ANA", 108, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Hoja 2: LP1912-215  (3 filas, sin filas nuevas)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"

$rows2 = @(
    @($newTime, "04:45", "215A_EL PATO", 2, "LP1912"),
    @($newTime, "05:34", "215B_EL PATO", 51, "LP1912"),
    @($newTime, "06:11", "215A_EL PATO", 88, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Hoja 3: 6203-6173  (2 filas -> 3 filas)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 3"

$rows3 = @(
    @($newTime, "05:43", "215A_LA PLATA", 60, "L6173"),
    @($newTime, "06:08", "215A_LA PLATA", 85, "L6173"),
    @($newTime, "06:32", "215C_LA PLATA", 109, "L6203")
)

$r = 6
foreach ($row in $rows3) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $r++
}
